# add: Agregar segunda historia de Seguro
#
# The document ends with an empty "Sin espaciado" paragraph that sits
# right before the final section break. This script inserts, immediately
# before that trailing separator paragraph, three new paragraphs that
# mirror the layout already used for "Historia de usuario no.1":
#
#   1. an empty "Sin espaciado" separator paragraph (same tab/spacing/
#      underline formatting as the other separators in the document),
#   2. the "Historia de usuario no.2" heading (sz=32 half-points, i.e.
#      16pt, matching the existing "Historia de usuario no.1" heading),
#   3. the user-story body paragraph (sz=24 half-points, i.e. 12pt,
#      matching the existing body paragraphs).
#
# We build the exact target OOXML for the three new paragraphs and hand
# it to Range.InsertXML so the inserted markup matches the surrounding
# document precisely (no inherited tab stops / spacing / underline
# leaking into the new heading & body paragraphs, and no stray pStyle on
# paragraphs that should fall back to the document's default "Normal"
# style).

$d = $word.ActiveDocument

# The very last paragraph in the main body is the empty "Sin espaciado"
# paragraph that closes out the document (right before the sectPr).
$lastPara = $d.Paragraphs.Last
$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(1)  # wdCollapseStart - insert right before it

$newStoryXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Sinespaciado"/>
              <w:tabs>
                <w:tab w:val="left" w:pos="5670"/>
              </w:tabs>
              <w:spacing w:line="276" w:lineRule="auto"/>
              <w:rPr>
                <w:u w:val="single"/>
              </w:rPr>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>Historia de usuario no.2: Actualizar datos de u seguro vigente.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>Como el encargado del área de seguros, quiero poder actualizar los datos de un seguro que ya está registrado, para poder reflejar cualquier cambio que haya surgido, como renovaciones de póliza, actualizaciones de cobertura o corrección de errores en la información. Necesito poder editar los campos necesarios, guardar cambios, y que el sistema mantenga un historial básico de las modificaciones realizadas. También quiero poder desactivar seguros vencidos o que ya no estén vigentes, sin tener que eliminarlos por completo.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($newStoryXml)
